$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Asistencia Q1 Agosto")
$src.Copy([System.Reflection.Missing]::Value, $src)
$newSheet = $wb.Worksheets.Item($src.Index + 1)
Write-Host "New sheet name:" $newSheet.Name

# Delete rows 17-31 (the "Arely siguenza" block)
$newSheet.Range("A17:C31").EntireRow.Delete()

# Change A2:A16 text to "Ana Flores"
$newSheet.Range("A2:A16").Value = "Ana Flores"

# Set scroll position on "Asistencia Q2 julio" sheet
$ws2 = $wb.Worksheets.Item("Asistencia Q2 julio")
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 16

# Update selection on new sheet and activate it last so it stays the visible tab
$newSheet.Activate()
$newSheet.Range("A2:A16").Select()

Write-Host "Dimension check done"
